# Add two new worksheets - "Debts" and "Fixed Assets" - after the existing
# "Sally" sheet, matching the structure used by "John"/"Sally" example
# workbooks (a bold header row, with a couple of currency-formatted columns).

$wb = $excel.ActiveWorkbook
$sally = $wb.Worksheets.Item(2)

# Create the sheets in final left-to-right order (Debts, then Fixed Assets)
# right after "Sally" so they get sheetId 3 / 4 respectively.
$debts = $wb.Worksheets.Add($null, $sally)
$debts.Name = "Debts"

$fa = $wb.Worksheets.Add($null, $debts)
$fa.Name = "Fixed Assets"

# Populate "Fixed Assets" first so its unique header strings land earlier in
# the shared-string table than "Debts"'s unique strings ("term"/"amount").
$fa.Range("A1").Value = "name"
$fa.Range("B1").Value = "type"
$fa.Range("C1").Value = "basis"
$fa.Range("D1").Value = "value"
$fa.Range("E1").Value = "rate"
$fa.Range("F1").Value = "yod"
$fa.Range("G1").Value = "commission"
$fa.Range("A1:G1").Font.Bold = $true
$fa.Range("C1:D1").NumberFormat = "`"$`"#,##0"
$null = $fa.Rows(1).Select()

# Now populate "Debts".
$debts.Range("A1").Value = "name"
$debts.Range("B1").Value = "type"
$debts.Range("C1").Value = "year"
$debts.Range("D1").Value = "term"
$debts.Range("E1").Value = "amount"
$debts.Range("F1").Value = "rate"
$debts.Range("A1:F1").Font.Bold = $true
$debts.Range("E1").NumberFormat = "`"$`"#,##0"
$null = $debts.Rows(1).Select()

# "Debts" is the tab that ends up active/selected in the saved workbook.
$debts.Activate()
